$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, pushing the existing "Test Co." row down to row 4
$ws.Rows("3:3").Insert()

# Copy the formatting of the row that shifted down (row 4, col B) up onto the
# newly inserted row 3 col B, so the new cell keeps the same "Comma" cell style.
$ws.Range("B4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's values
$ws.Range("A3").Value = 101
$ws.Range("B3").Value = "Amped Well Servicing Ltd."

# Update the selection to match where the user left off
$ws.Range("B4").Select() | Out-Null
